# Update the cryptocurrency price/volume table to the latest scrape.
# A leading apostrophe forces text entry (mirrors typing into Excel),
# so numeric-looking strings like "1.008" or "27.000.76" stay text
# instead of being auto-coerced to numbers/dates, matching the source cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.000.76"
$ws.Range("E2").Value = "'  -1.23%  "

$ws.Range("D3").Value = "'1.822.44"
$ws.Range("E3").Value = "'  -0.61%  "

$ws.Range("D4").Value = "'1.009"
$ws.Range("E4").Value = "'  -0.28%  "

$ws.Range("D5").Value = "'309.63"
$ws.Range("E5").Value = "'  -1.66%  "

$ws.Range("D6").Value = "'1.007"

$ws.Range("D7").Value = "'0.4649"
$ws.Range("E7").Value = "'  -2.05%  "

$ws.Range("D8").Value = "'0.3663"
$ws.Range("E8").Value = "'  -0.76%  "

$ws.Range("D9").Value = "'0.07235"
$ws.Range("E9").Value = "'  -2.99%  "

$ws.Range("D10").Value = "'0.8598"
$ws.Range("E10").Value = "'  -3.02%  "

$ws.Range("D11").Value = "'19.84"
$ws.Range("E11").Value = "'  -3.05%  "

$ws.Range("B12").Value = "'WrappedEther"
$ws.Range("C12").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.887.43"
$ws.Range("E12").Value = "'  +0.43%  "

$ws.Range("B13").Value = "'TRON"
$ws.Range("C13").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.07599"
$ws.Range("E13").Value = "'  +3.58%  "

$ws.Range("D14").Value = "'5.335"
$ws.Range("E14").Value = "'  -1.94%  "

$ws.Range("D15").Value = "'6.496"
$ws.Range("E15").Value = "'  -1.31%  "

$ws.Range("D16").Value = "'91.62"
$ws.Range("E16").Value = "'  -1.86%  "

$ws.Range("D17").Value = "'1.009"
$ws.Range("E17").Value = "'  -0.07%  "

$ws.Range("D18").Value = "'0.000008634"
$ws.Range("E18").Value = "'  -1.87%  "

$ws.Range("D19").Value = "'1.006"
$ws.Range("E19").Value = "'  -0.41%  "

$ws.Range("B20").Value = "'Avalanche"
$ws.Range("C20").Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'14.49"

$ws.Range("B21").Value = "'WrappedBTC"
$ws.Range("C21").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").Value = "'26.816.20"
$ws.Range("E21").Value = "'  -2.77%  "

$ws.Range("D22").Value = "'5.153"
$ws.Range("E22").Value = "'  -2.61%  "

$ws.Range("D23").Value = "'10.53"
$ws.Range("E23").Value = "'  -1.25%  "

$ws.Range("D24").Value = "'2.046.83"
$ws.Range("E24").Value = "'  -2.19%  "

$ws.Range("D25").Value = "'151.53"
$ws.Range("E25").Value = "'  -0.31%  "

$ws.Range("D26").Value = "'1.842"
$ws.Range("E26").Value = "'  -2.65%  "

$ws.Range("D27").Value = "'18.15"
$ws.Range("E27").Value = "'  -2.73%  "

$ws.Range("D28").Value = "'2.052"
$ws.Range("E28").Value = "'  -4.43%  "

$ws.Range("E29").Value = "'  -2.56%  "

$ws.Range("D30").Value = "'115.25"
$ws.Range("E30").Value = "'  -1.75%  "

$ws.Range("D31").Value = "'0.08840"
$ws.Range("E31").Value = "'  -1.71%  "

$ws.Range("D32").Value = "'2.955"
$ws.Range("E32").Value = "'  +0.29%  "

$ws.Range("D33").Value = "'4.428"
$ws.Range("E33").Value = "'  -2.65%  "

$ws.Range("D34").Value = "'1.128"
$ws.Range("E34").Value = "'  -4.13%  "

$ws.Range("D35").Value = "'0.7195"
$ws.Range("E35").Value = "'  -4.46%  "

$ws.Range("D36").Value = "'1.077"
$ws.Range("E36").Value = "'  -2.33%  "

$ws.Range("D37").Value = "'0.05254"
$ws.Range("E37").Value = "'  -1.72%  "

$ws.Range("E38").Value = "'  -1.61%  "

$ws.Range("D39").Value = "'2.402"
$ws.Range("E39").Value = "'  +0.52%  "

$ws.Range("D40").Value = "'2.930"
$ws.Range("E40").Value = "'  -1.60%  "

$ws.Range("D41").Value = "'7.146"
$ws.Range("E41").Value = "'  -1.80%  "

$ws.Range("D42").Value = "'0.5159"
$ws.Range("E42").Value = "'  -2.95%  "

$ws.Range("B43").Value = "'Frax"
$ws.Range("C43").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").Value = "'0.8718"
$ws.Range("E43").Value = "'  -13.75%  "

$ws.Range("B44").Value = "'Algorand"
$ws.Range("C44").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "'0.1628"
$ws.Range("E44").Value = "'  -1.96%  "

$ws.Range("D45").Value = "'8.167"
$ws.Range("E45").Value = "'  -3.83%  "

$ws.Range("D46").Value = "'0.4809"
$ws.Range("E46").Value = "'  -2.12%  "

$ws.Range("E47").Value = "'  -0.29%  "

$ws.Range("D48").Value = "'10.19"
$ws.Range("E48").Value = "'  -3.54%  "

$ws.Range("D49").Value = "'102.65"
$ws.Range("E49").Value = "'  -2.39%  "

$ws.Range("D50").Value = "'1.621"
$ws.Range("E50").Value = "'  -3.07%  "

$ws.Range("D51").Value = "'0.06244"
